$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows to match repulled data.
$updates = @{
    6  = -7
    7  = -5
    11 = -1
    12 = 6
    16 = 2
    23 = 3
    24 = 4
    32 = 3
    35 = -1
    36 = -2
    41 = 0
    51 = 0
    57 = 4
    58 = -3
    59 = -4
    60 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
